$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "27.940.16"
$ws.Range("E2").Value = "  +1.55%  "
$ws.Range("D3").Value = "1.641.25"
$ws.Range("E3").Value = "  +1.21%  "
$ws.Range("E4").Value = "  +0.01%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "213.56"
$ws.Range("D5").ClearFormats()
$ws.Range("E5").Value = "  +0.96%  "
$ws.Range("E6").Value = "  +0.15%  "
$ws.Range("E7").Value = "  -0.03%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "23.81"
$ws.Range("D8").ClearFormats()
$ws.Range("E8").Value = "  +3.29%  "
$ws.Range("E9").Value = "  +0.21%  "
$ws.Range("E10").Value = "  +0.84%  "
$ws.Range("E11").Value = "  -0.61%  "
$ws.Range("D12").Value = "1.873.94"
$ws.Range("E12").Value = "  +1.22%  "
$ws.Range("D13").Value = "1.641.49"
$ws.Range("E13").Value = "  +0.98%  "
$ws.Range("B14").Value = "Polygon"
$ws.Range("C14").Value = "https://coinranking.com/coin/uW2tk-ILY0ii+polygon-matic"
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "0.577"
$ws.Range("D14").ClearFormats()
$ws.Range("E14").Value = "  +5.02%  "
$ws.Range("B15").Value = "Polkadot"
$ws.Range("C15").Value = "https://coinranking.com/coin/25W7FG7om+polkadot-dot"
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "4.09"
$ws.Range("D15").ClearFormats()
$ws.Range("E15").Value = "  +1.16%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "66.18"
$ws.Range("D16").ClearFormats()
$ws.Range("D17").Value = "27.929.94"
$ws.Range("E17").Value = "  +1.59%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "231.57"
$ws.Range("D18").ClearFormats()
$ws.Range("E18").Value = "  +0.72%  "
$ws.Range("E19").Value = "  +1.18%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "7.61"
$ws.Range("D20").ClearFormats()
$ws.Range("E20").Value = "  +0.85%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "11.62"
$ws.Range("D21").ClearFormats()
$ws.Range("E21").Value = "  +11.97%  "
$ws.Range("E22").Value = "  -0.05%  "
$ws.Range("E23").Value = "  +1.45%  "
$ws.Range("E24").Value = "  -2.21%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "152.07"
$ws.Range("D25").ClearFormats()
$ws.Range("E25").Value = "  +1.69%  "
$ws.Range("E26").Value = "  +0.98%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "15.76"
$ws.Range("D27").ClearFormats()
$ws.Range("E27").Value = "  +1.46%  "
$ws.Range("E28").Value = "  +0.83%  "
$ws.Range("E29").Value = "  -0.05%  "
$ws.Range("E30").Value = "  +1.01%  "
$ws.Range("E31").Value = "  +0.48%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "3.34"
$ws.Range("D32").ClearFormats()
$ws.Range("E32").Value = "  +2.11%  "
$ws.Range("D33").Value = "1.425.01"
$ws.Range("E33").Value = "  -2.68%  "
$ws.Range("E34").Value = "  +2.27%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "1.58"
$ws.Range("D35").ClearFormats()
$ws.Range("E35").Value = "  +1.96%  "
$ws.Range("E36").Value = "  +0.39%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.894"
$ws.Range("D37").ClearFormats()
$ws.Range("E37").Value = "  +2.54%  "
$ws.Range("E38").Value = "  +0.71%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.924"
$ws.Range("D39").ClearFormats()
$ws.Range("E39").Value = "  -1.97%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.557"
$ws.Range("D40").ClearFormats()
$ws.Range("E40").Value = "  +1.04%  "
$ws.Range("E41").Value = "  +2.27%  "
$ws.Range("E42").Value = "  -0.06%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "67.34"
$ws.Range("D43").ClearFormats()
$ws.Range("E43").Value = "  +0.00%  "
$ws.Range("E44").Value = "  +0.68%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "1.83"
$ws.Range("D45").ClearFormats()
$ws.Range("E45").Value = "  +4.58%  "
$ws.Range("E46").Value = "  +2.82%  "
$ws.Range("E47").Value = "  +0.22%  "
$ws.Range("D48").Value = "1.783.05"
$ws.Range("E48").Value = "  +1.26%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "88.72"
$ws.Range("D49").ClearFormats()
$ws.Range("E49").Value = "  +1.76%  "
$ws.Range("E51").Value = "  +0.69%  "
